$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = 2.180165333333334
$ws.Range("H2").Value = 6.540496
$ws.Range("I2").Value = 0.01970539991828544
$ws.Range("J2").Value = 0.01970539991828544
$ws.Range("M2").Value = 227.11144
$ws.Range("N2").Value = 681.33432
$ws.Range("O2").Value = 0.8625743548356182
$ws.Range("P2").Value = 0.8625743548356182
$ws.Range("Q2").Value = 495.1404882914134
$ws.Range("R2").Value = 4456.26439462272
$ws.Range("S2").Value = 0.0169973726212929
$ws.Range("T2").Value = 0.0169973726212929
$ws.Range("G3").Value = 2.180165333333334
$ws.Range("H3").Value = 6.540496
$ws.Range("I3").Value = 0.01970539991828544
$ws.Range("J3").Value = 0.01970539991828544
$ws.Range("O3").Value = 0.001598666154760757
$ws.Range("P3").Value = 0.001598666154760757
$ws.Range("Q3").Value = 0.9176766455502224
$ws.Range("R3").Value = 8.259089809952002
$ws.Range("S3").Value = 0.00003150235591538831
$ws.Range("T3").Value = 0.0000315023559153883
$ws.Range("G4").Value = 2.180165333333334
$ws.Range("H4").Value = 6.540496
$ws.Range("I4").Value = 0.01970539991828544
$ws.Range("J4").Value = 0.01970539991828544
$ws.Range("M4").Value = 3.233093
$ws.Range("N4").Value = 9.699279000000001
$ws.Range("O4").Value = 0.01227935989749593
$ws.Range("P4").Value = 0.01227935989749593
$ws.Range("Q4").Value = 7.048677278042668
$ws.Range("R4").Value = 63.438095502384
$ws.Range("S4").Value = 0.0002419696975207138
$ws.Range("T4").Value = 0.0002419696975207138
$ws.Range("G5").Value = 2.180165333333334
$ws.Range("H5").Value = 6.540496
$ws.Range("I5").Value = 0.01970539991828544
$ws.Range("J5").Value = 0.01970539991828544
$ws.Range("M5").Value = 32.52945966666667
$ws.Range("N5").Value = 97.588379
$ws.Range("O5").Value = 0.1235476191121251
$ws.Range("P5").Value = 0.1235476191121251
$ws.Range("Q5").Value = 70.91960027733157
$ws.Range("R5").Value = 638.276402495984
$ws.Range("S5").Value = 0.002434555243556431
$ws.Range("T5").Value = 0.002434555243556431
$ws.Range("I6").Value = 0.733713204346044
$ws.Range("J6").Value = 0.7337132043460441
$ws.Range("M6").Value = 227.11144
$ws.Range("N6").Value = 681.33432
$ws.Range("O6").Value = 0.8625743548356182
$ws.Range("P6").Value = 0.8625743548356182
$ws.Range("Q6").Value = 18436.11983376422
$ws.Range("R6").Value = 165925.0785038779
$ws.Range("S6").Value = 0.632882193873163
$ws.Range("T6").Value = 0.6328821938731631
$ws.Range("I7").Value = 0.733713204346044
$ws.Range("J7").Value = 0.7337132043460441
$ws.Range("O7").Value = 0.001598666154760757
$ws.Range("P7").Value = 0.001598666154760757
$ws.Range("S7").Value = 0.001172962467089084
$ws.Range("T7").Value = 0.001172962467089084
$ws.Range("I8").Value = 0.733713204346044
$ws.Range("J8").Value = 0.7337132043460441
$ws.Range("M8").Value = 3.233093
$ws.Range("N8").Value = 9.699279000000001
$ws.Range("O8").Value = 0.01227935989749593
$ws.Range("P8").Value = 0.01227935989749593
$ws.Range("Q8").Value = 262.4512881504527
$ws.Range("R8").Value = 2362.061593354074
$ws.Range("S8").Value = 0.009009528497710051
$ws.Range("T8").Value = 0.009009528497710051
$ws.Range("I9").Value = 0.733713204346044
$ws.Range("J9").Value = 0.7337132043460441
$ws.Range("M9").Value = 32.52945966666667
$ws.Range("N9").Value = 97.588379
$ws.Range("O9").Value = 0.1235476191121251
$ws.Range("P9").Value = 0.1235476191121251
$ws.Range("Q9").Value = 2640.628832005408
$ws.Range("R9").Value = 23765.65948804867
$ws.Range("S9").Value = 0.09064851950808189
$ws.Range("T9").Value = 0.09064851950808191
$ws.Range("G10").Value = 25.672264
$ws.Range("H10").Value = 77.016792
$ws.Range("I10").Value = 0.2320384702908474
$ws.Range("J10").Value = 0.2320384702908474
$ws.Range("M10").Value = 227.11144
$ws.Range("N10").Value = 681.33432
$ws.Range("O10").Value = 0.8625743548356182
$ws.Range("P10").Value = 0.8625743548356182
$ws.Range("Q10").Value = 5830.46484510016
$ws.Range("R10").Value = 52474.18360590144
$ws.Range("S10").Value = 0.2001504338081715
$ws.Range("T10").Value = 0.2001504338081715
$ws.Range("G11").Value = 25.672264
$ws.Range("H11").Value = 77.016792
$ws.Range("I11").Value = 0.2320384702908474
$ws.Range("J11").Value = 0.2320384702908474
$ws.Range("O11").Value = 0.001598666154760757
$ws.Range("P11").Value = 0.001598666154760757
$ws.Range("Q11").Value = 10.80598647772267
$ws.Range("R11").Value = 97.25387829950401
$ws.Range("S11").Value = 0.0003709520490564371
$ws.Range("T11").Value = 0.000370952049056437
$ws.Range("G12").Value = 25.672264
$ws.Range("H12").Value = 77.016792
$ws.Range("I12").Value = 0.2320384702908474
$ws.Range("J12").Value = 0.2320384702908474
$ws.Range("M12").Value = 3.233093
$ws.Range("N12").Value = 9.699279000000001
$ws.Range("O12").Value = 0.01227935989749593
$ws.Range("P12").Value = 0.01227935989749593
$ws.Range("Q12").Value = 83.00081703255201
$ws.Range("R12").Value = 747.007353292968
$ws.Range("S12").Value = 0.002849283886765733
$ws.Range("T12").Value = 0.002849283886765732
$ws.Range("G13").Value = 25.672264
$ws.Range("H13").Value = 77.016792
$ws.Range("I13").Value = 0.2320384702908474
$ws.Range("J13").Value = 0.2320384702908474
$ws.Range("M13").Value = 32.52945966666667
$ws.Range("N13").Value = 97.588379
$ws.Range("O13").Value = 0.1235476191121251
$ws.Range("P13").Value = 0.1235476191121251
$ws.Range("Q13").Value = 835.1048763400187
$ws.Range("R13").Value = 7515.943887060168
$ws.Range("S13").Value = 0.02866780054685378
$ws.Range("T13").Value = 0.02866780054685378
$ws.Range("G14").Value = 1.608999666666667
$ws.Range("H14").Value = 4.826999
$ws.Range("I14").Value = 0.01454292544482312
$ws.Range("J14").Value = 0.01454292544482312
$ws.Range("M14").Value = 227.11144
$ws.Range("N14").Value = 681.33432
$ws.Range("O14").Value = 0.8625743548356182
$ws.Range("P14").Value = 0.8625743548356182
$ws.Range("Q14").Value = 365.4222312561867
$ws.Range("R14").Value = 3288.80008130568
$ws.Range("S14").Value = 0.0125443545329908
$ws.Range("T14").Value = 0.0125443545329908
$ws.Range("G15").Value = 1.608999666666667
$ws.Range("H15").Value = 4.826999
$ws.Range("I15").Value = 0.01454292544482312
$ws.Range("J15").Value = 0.01454292544482312
$ws.Range("O15").Value = 0.001598666154760757
$ws.Range("P15").Value = 0.001598666154760757
$ws.Range("Q15").Value = 0.6772612123597779
$ws.Range("R15").Value = 6.095350911238
$ws.Range("S15").Value = 0.00002324928269984775
$ws.Range("T15").Value = 0.00002324928269984775
$ws.Range("G16").Value = 1.608999666666667
$ws.Range("H16").Value = 4.826999
$ws.Range("I16").Value = 0.01454292544482312
$ws.Range("J16").Value = 0.01454292544482312
$ws.Range("M16").Value = 3.233093
$ws.Range("N16").Value = 9.699279000000001
$ws.Range("O16").Value = 0.01227935989749593
$ws.Range("P16").Value = 0.01227935989749593
$ws.Range("Q16").Value = 5.202045559302333
$ws.Range("R16").Value = 46.818410033721
$ws.Range("S16").Value = 0.0001785778154994343
$ws.Range("T16").Value = 0.0001785778154994342
$ws.Range("G17").Value = 1.608999666666667
$ws.Range("H17").Value = 4.826999
$ws.Range("I17").Value = 0.01454292544482312
$ws.Range("J17").Value = 0.01454292544482312
$ws.Range("M17").Value = 32.52945966666667
$ws.Range("N17").Value = 97.588379
$ws.Range("O17").Value = 0.1235476191121251
$ws.Range("P17").Value = 0.1235476191121251
$ws.Range("Q17").Value = 52.33988976051344
$ws.Range("R17").Value = 471.059007844621
$ws.Range("S17").Value = 0.00179674381363304
$ws.Range("T17").Value = 0.00179674381363304
